# Increase confidence threshold and update signal generation parameters
# Applies the regenerated genx_signals.xlsx snapshot to the workbook:
#   - Active Signals: only 2 signals remain active now (rows 4-9 removed),
#     and the two surviving rows carry fresh values.
#   - Summary Dashboard: refresh the aggregate counters.
#   - Signal History: the rolling history window advanced, so every row's
#     contents shift to the newer batch of signals.
#   - The SELL/BUY highlight colors in the shared palette are swapped
#     (SELL now reads as green, BUY now reads as red).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Active Signals"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Active Signals")

# Remove the six signals that dropped off the active list (old rows 4-9).
# Deleting from the bottom up keeps the remaining row numbers stable.
$ws1.Rows.Item(9).Delete()
$ws1.Rows.Item(8).Delete()
$ws1.Rows.Item(7).Delete()
$ws1.Rows.Item(6).Delete()
$ws1.Rows.Item(5).Delete()
$ws1.Rows.Item(4).Delete()

# Row 2 now reports the EURUSD BUY signal.
$ws1.Cells.Item(2, 1).Value = "2025-07-28 20:18"
$ws1.Cells.Item(2, 2).Value = "EURUSD"
$ws1.Cells.Item(2, 3).Value = "BUY"
$ws1.Cells.Item(2, 4).Value = 1.10095
$ws1.Cells.Item(2, 5).Value = 1.09751
$ws1.Cells.Item(2, 6).Value = 1.10812
$ws1.Cells.Item(2, 7).Value = 0.09
$ws1.Cells.Item(2, 8).Value = "76.0%"
$ws1.Cells.Item(2, 9).Value = 2.09
$ws1.Cells.Item(2, 10).Value = "Active"

# Row 3 now reports the XAUUSD SELL signal.
$ws1.Cells.Item(3, 1).Value = "2025-07-28 20:04"
$ws1.Cells.Item(3, 2).Value = "XAUUSD"
$ws1.Cells.Item(3, 3).Value = "SELL"
$ws1.Cells.Item(3, 4).Value = 2646.76589
$ws1.Cells.Item(3, 5).Value = 2646.76804
$ws1.Cells.Item(3, 6).Value = 2646.75787
$ws1.Cells.Item(3, 7).Value = 0.07000000000000001
$ws1.Cells.Item(3, 8).Value = "76.0%"
$ws1.Cells.Item(3, 9).Value = 3.73
$ws1.Cells.Item(3, 10).Value = "Active"

# The SELL/BUY conditional colors were swapped in the shared fill palette
# (SELL: was pink FFC7CE -> now green C6EFCE; BUY: was green C6EFCE -> now
# pink FFC7CE). Re-apply so the remaining signal rows pick up the new look.
$ws1.Cells.Item(2, 3).Interior.Color = 13551615  # BUY  -> FFC7CE
$ws1.Cells.Item(3, 3).Interior.Color = 13561798  # SELL -> C6EFCE

# ---------------------------------------------------------------------------
# Sheet 2: "Summary Dashboard"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary Dashboard")

$ws2.Cells.Item(4, 2).Value = 2          # Active Signals
$ws2.Cells.Item(5, 2).Value = 8          # BUY Signals
$ws2.Cells.Item(6, 2).Value = 7          # SELL Signals
$ws2.Cells.Item(7, 2).Value = "84.9%"    # Average Confidence
$ws2.Cells.Item(8, 2).Value = "2.29"     # Average Risk/Reward
$ws2.Cells.Item(9, 2).Value = "2025-07-28 19:56:12"  # Last Update

# ---------------------------------------------------------------------------
# Sheet 3: "Signal History"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Signal History")

function Set-HistoryRow {
    param($sheet, $row, $timestamp, $symbol, $signal, $entry, $sl, $tp, $lots, $confidence, $rr, $status)

    $sheet.Cells.Item($row, 1).Value = $timestamp
    $sheet.Cells.Item($row, 2).Value = $symbol
    $sheet.Cells.Item($row, 3).Value = $signal
    $sheet.Cells.Item($row, 4).Value = $entry
    $sheet.Cells.Item($row, 5).Value = $sl
    $sheet.Cells.Item($row, 6).Value = $tp
    $sheet.Cells.Item($row, 7).Value = $lots
    $sheet.Cells.Item($row, 8).Value = $confidence
    $sheet.Cells.Item($row, 9).Value = $rr
    $sheet.Cells.Item($row, 10).Value = $status
}

Set-HistoryRow $ws3 2  "2025-07-28 19:57" "EURUSD" "BUY"  1.10417     1.09996     1.11413     0.09 0.85 2.36 "Pending"
Set-HistoryRow $ws3 3  "2025-07-28 19:27" "USDJPY" "SELL" 150.17889   150.67464   149.39819   0.05 0.88 1.57 "Filled"
Set-HistoryRow $ws3 4  "2025-07-28 20:06" "XAUGBP" "SELL" 2105.81156  2105.81421  2105.80399  0.01 0.85 2.85 "Pending"
Set-HistoryRow $ws3 5  "2025-07-28 20:00" "USDJPY" "BUY"  150.23958   150.00475   150.65111   0.06 0.87 1.75 "Pending"
Set-HistoryRow $ws3 6  "2025-07-28 19:52" "XAUCHF" "BUY"  2345.83131  2345.82699  2345.83617  0.07000000000000001 0.88 1.12 "Filled"
Set-HistoryRow $ws3 7  "2025-07-28 19:34" "EURUSD" "SELL" 1.10132     1.10507     1.09573     0.05 0.8  1.49 "Filled"
Set-HistoryRow $ws3 8  "2025-07-28 20:18" "EURUSD" "BUY"  1.10095     1.09751     1.10812     0.09 0.76 2.09 "Active"
Set-HistoryRow $ws3 9  "2025-07-28 19:50" "XAUGBP" "SELL" 2093.56937  2093.57231  2093.5647   0.02 0.79 1.59 "Pending"
Set-HistoryRow $ws3 10 "2025-07-28 19:42" "NZDUSD" "BUY"  0.59234     0.58994     0.60129     0.07000000000000001 0.92 3.74 "Filled"
Set-HistoryRow $ws3 11 "2025-07-28 19:30" "NZDUSD" "BUY"  0.5891999999999999 0.58645 0.59809   0.08 0.93 3.24 "Pending"
Set-HistoryRow $ws3 12 "2025-07-28 20:04" "XAUCAD" "BUY"  3602.70636  3602.70329  3602.71461  0.07000000000000001 0.8100000000000001 2.68 "Filled"
Set-HistoryRow $ws3 13 "2025-07-28 20:03" "USDCHF" "SELL" 0.87935     0.8829900000000001 0.87307 0.1 0.9399999999999999 1.72 "Pending"
Set-HistoryRow $ws3 14 "2025-07-28 19:30" "XAUCHF" "SELL" 2341.82489  2341.82898  2341.81819  0.01 0.76 1.64 "Filled"
Set-HistoryRow $ws3 15 "2025-07-28 20:13" "USDCHF" "BUY"  0.88243     0.8788899999999999 0.89208 0.08 0.9399999999999999 2.73 "Filled"
Set-HistoryRow $ws3 16 "2025-07-28 20:04" "XAUUSD" "SELL" 2646.76589  2646.76804  2646.75787  0.07000000000000001 0.76 3.73 "Active"
